# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-order Panama/Portugal/Singapur block in the "Pais" column (A41:A45) ---
# Current order: Bolivia, Portugal, Singapur, Panama, Republica Dominicana
# Target order : Bolivia, Panama, Portugal, Singapur, Republica Dominicana
$ws.Range("A42").Value = "Panama"
$ws.Range("A43").Value = "Portugal"
$ws.Range("A44").Value = "Singapur"

# --- Update the "last updated" timestamp string (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 14 de Julio de 2020 a las 02:13"

# --- Refresh per-country statistics (Casos totales, Nuevos casos, Casos activos,
#     Recuperados, Casos criticos, Muertes hoy, Muertes) ---

# Row 4 (Estados Unidos)
$ws.Range("B4").Value = 3478481
$ws.Range("C4").Value = 64486
$ws.Range("D4").Value = 1549112
$ws.Range("E4").Value = 1791135
$ws.Range("G4").Value = 452
$ws.Range("H4").Value = 138234

# Row 5 (Brasil)
$ws.Range("B5").Value = 1887959
$ws.Range("C5").Value = 21783
$ws.Range("E5").Value = 601526
$ws.Range("G5").Value = 770
$ws.Range("H5").Value = 72921

# Row 23 (Canada)
$ws.Range("B23").Value = 108155
$ws.Range("C23").Value = 565
$ws.Range("D23").Value = 71841
$ws.Range("E23").Value = 27524
$ws.Range("G23").Value = 7
$ws.Range("H23").Value = 8790

# Row 42 (now Panama)
$ws.Range("B42").Value = 47173
$ws.Range("C42").Value = 1540
$ws.Range("D42").Value = 23919
$ws.Range("E42").Value = 22322
$ws.Range("G42").Value = 23
$ws.Range("H42").Value = 932

# Row 43 (now Portugal)
$ws.Range("B43").Value = 46818
$ws.Range("C43").Value = 306
$ws.Range("D43").Value = 31065
$ws.Range("E43").Value = 14091
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 1662

# Row 44 (now Singapur)
$ws.Range("B44").Value = 46283
$ws.Range("C44").Value = 322
$ws.Range("D44").Value = 42541
$ws.Range("E44").Value = 3716
$ws.Range("H44").Value = 26

# Row 134
$ws.Range("B134").Value = 1287
$ws.Range("C134").Value = 66
$ws.Range("D134").Value = 330
$ws.Range("E134").Value = 933
$ws.Range("G134").Value = 1
$ws.Range("H134").Value = 24

# Row 145
$ws.Range("B145").Value = 989
$ws.Range("C145").Value = 2
$ws.Range("D145").Value = 903
$ws.Range("E145").Value = 55

# Row 148
$ws.Range("B148").Value = 858
$ws.Range("C148").Value = 3
$ws.Range("E148").Value = 3

# Row 151
$ws.Range("B151").Value = 732
$ws.Range("C151").Value = 3
$ws.Range("E151").Value = 432

# Row 152
$ws.Range("B152").Value = 721
$ws.Range("C152").Value = 1
$ws.Range("D152").Value = 517
$ws.Range("E152").Value = 189

# Row 167
$ws.Range("B167").Value = 300
$ws.Range("C167").Value = 3
$ws.Range("D167").Value = 155
$ws.Range("E167").Value = 128

# Row 173
$ws.Range("B173").Value = 203
$ws.Range("C173").Value = 2
$ws.Range("D173").Value = 200
$ws.Range("E173").Value = 2

# Row 176
$ws.Range("D176").Value = 179
$ws.Range("E176").Value = 1

# Row 180
$ws.Range("D180").Value = 124
$ws.Range("E180").Value = 1
